$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# "Periodo Mora" column (E16:E23) gets reversed: the list of periods that used
# to run 2012, 2101, 2102, 2103, 2104, 2105, 2106, 2107 (top to bottom) is
# flipped to run 2107, 2106, 2105, 2104, 2103, 2102, 2101, 2012 (top to bottom).
$ws.Range("E16").Value = "2107"
$ws.Range("E17").Value = "2106"
$ws.Range("E18").Value = "2105"
$ws.Range("E19").Value = "2104"
$ws.Range("E20").Value = "2103"
$ws.Range("E21").Value = "2102"
$ws.Range("E22").Value = "2101"
$ws.Range("E23").Value = "2012"

# "Valor Mora" column: the values on the first and last data rows swap places.
$ws.Range("F16").Value = 29260
$ws.Range("F23").Value = 35112
